$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (values chosen so the saved OOXML <col width> lands
# on the exact target after Excel's character-width/pixel rounding)
$ws.Columns.Item(1).ColumnWidth = 11.15
$ws.Columns.Item(2).ColumnWidth = 11.15
$ws.Columns.Item(3).ColumnWidth = 13.15

# Update row 2 data
$ws.Range("A2").Value = "GEOJENKINS"
$ws.Range("B2").Value = "geojenkins"
$ws.Range("C2").Value = "10.181.4.105"
$ws.Range("H2").Value = "No instalado"
$ws.Range("I2").Value = "1.2.0.956"
$ws.Range("J2").Value = "active"
$ws.Range("K2").Value = "inactive"
